$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 6.479741854214802
$ws.Cells.Item(2, 4).Value2 = 10.1193391787791
$ws.Cells.Item(2, 5).Value2 = 14.05621295591393
$ws.Cells.Item(2, 6).Value2 = 28.39693932865039
$ws.Cells.Item(2, 7).Value2 = 26.93368770125346
$ws.Cells.Item(2, 8).Value2 = 13.89643175932336
$ws.Cells.Item(2, 9).Value2 = 19.16038114729668
$ws.Cells.Item(2, 10).Value2 = 9.754796521151899
$ws.Cells.Item(2, 11).Value2 = 11.60143225170288
$ws.Cells.Item(2, 15).Value2 = 20.8971927117789
$ws.Cells.Item(3, 2).Value2 = 6.312005066377872
$ws.Cells.Item(3, 4).Value2 = 10.05046149670758
$ws.Cells.Item(3, 5).Value2 = 13.98639667277381
$ws.Cells.Item(3, 6).Value2 = 28.46183746384721
$ws.Cells.Item(3, 7).Value2 = 27.03561526224394
$ws.Cells.Item(3, 8).Value2 = 13.95132682990074
$ws.Cells.Item(3, 9).Value2 = 19.27309509879225
$ws.Cells.Item(3, 10).Value2 = 9.760658313036693
$ws.Cells.Item(3, 11).Value2 = 11.10622106096201
$ws.Cells.Item(3, 15).Value2 = 20.98896411666256
$ws.Cells.Item(4, 2).Value2 = 6.207256546673343
$ws.Cells.Item(4, 4).Value2 = 10.00957581456908
$ws.Cells.Item(4, 5).Value2 = 13.94622305105211
$ws.Cells.Item(4, 6).Value2 = 28.50945971369497
$ws.Cells.Item(4, 7).Value2 = 27.10842452916109
$ws.Cells.Item(4, 8).Value2 = 13.98748116101448
$ws.Cells.Item(4, 9).Value2 = 19.34626150236928
$ws.Cells.Item(4, 10).Value2 = 9.765818234251771
$ws.Cells.Item(4, 11).Value2 = 10.79030290331223
$ws.Cells.Item(4, 15).Value2 = 21.05037608163059
$ws.Cells.Item(5, 2).Value2 = 6.164188109753275
$ws.Cells.Item(5, 4).Value2 = 9.993281924866796
$ws.Cells.Item(5, 5).Value2 = 13.93054263609502
$ws.Cells.Item(5, 6).Value2 = 28.53081561953762
$ws.Cells.Item(5, 7).Value2 = 27.14065098111383
$ws.Cells.Item(5, 8).Value2 = 14.00282993662217
$ws.Cells.Item(5, 9).Value2 = 19.37707480589347
$ws.Cells.Item(5, 10).Value2 = 9.768313810753511
$ws.Cells.Item(5, 11).Value2 = 10.65873631755993
$ws.Cells.Item(5, 15).Value2 = 21.07667224262833
$ws.Cells.Item(6, 2).Value2 = 6.157015322992955
$ws.Cells.Item(6, 4).Value2 = 9.990598920259306
$ws.Cells.Item(6, 5).Value2 = 13.9279809985142
$ws.Cells.Item(6, 6).Value2 = 28.53447932015997
$ws.Cells.Item(6, 7).Value2 = 27.14615608031772
$ws.Cells.Item(6, 8).Value2 = 14.00541576636616
$ws.Cells.Item(6, 9).Value2 = 19.38225162359731
$ws.Cells.Item(6, 10).Value2 = 9.76875193707197
$ws.Cells.Item(6, 11).Value2 = 10.63672394094123
$ws.Cells.Item(6, 15).Value2 = 21.0811153316111
$ws.Cells.Item(7, 2).Value2 = 6.206677178756009
$ws.Cells.Item(7, 4).Value2 = 10.00935456407184
$ws.Cells.Item(7, 5).Value2 = 13.94600876619784
$ws.Cells.Item(7, 6).Value2 = 28.50973984141439
$ws.Cells.Item(7, 7).Value2 = 27.10884881755814
$ws.Cells.Item(7, 8).Value2 = 13.98768566760519
$ws.Cells.Item(7, 9).Value2 = 19.34667302058856
$ws.Cells.Item(7, 10).Value2 = 9.765850299317266
$ws.Cells.Item(7, 11).Value2 = 10.78853977910278
$ws.Cells.Item(7, 15).Value2 = 21.05072558161999
$ws.Cells.Item(8, 2).Value2 = 6.422305525904645
$ws.Cells.Item(8, 4).Value2 = 10.09530697354266
$ws.Cells.Item(8, 5).Value2 = 14.03158943481443
$ws.Cells.Item(8, 6).Value2 = 28.41769946756887
$ws.Cells.Item(8, 7).Value2 = 26.96670120431344
$ws.Cells.Item(8, 8).Value2 = 13.91485131794706
$ws.Cells.Item(8, 9).Value2 = 19.19842414347034
$ws.Cells.Item(8, 10).Value2 = 9.756494003409149
$ws.Cells.Item(8, 11).Value2 = 11.43323022728564
$ws.Cells.Item(8, 15).Value2 = 20.92778262519445
$ws.Cells.Item(9, 2).Value2 = 6.828777154910239
$ws.Cells.Item(9, 4).Value2 = 10.2743742861344
$ws.Cells.Item(9, 5).Value2 = 14.2201620323443
$ws.Cells.Item(9, 6).Value2 = 28.29911640174824
$ws.Cells.Item(9, 7).Value2 = 26.76971954325784
$ws.Cells.Item(9, 8).Value2 = 13.79145636532852
$ws.Cells.Item(9, 9).Value2 = 18.93904770217095
$ws.Cells.Item(9, 10).Value2 = 9.750509088258834
$ws.Cells.Item(9, 11).Value2 = 12.59764311007743
$ws.Cells.Item(9, 15).Value2 = 20.72700859903971
$ws.Cells.Item(10, 2).Value2 = 7.114436374926236
$ws.Cells.Item(10, 4).Value2 = 10.41146080765787
$ws.Cells.Item(10, 5).Value2 = 14.37046074341832
$ws.Cells.Item(10, 6).Value2 = 28.24999568704956
$ws.Cells.Item(10, 7).Value2 = 26.67565033294592
$ws.Cells.Item(10, 8).Value2 = 13.71265167142893
$ws.Cells.Item(10, 9).Value2 = 18.76748039223007
$ws.Cells.Item(10, 10).Value2 = 9.753614970396535
$ws.Cells.Item(10, 11).Value2 = 13.3858549115165
$ws.Cells.Item(10, 15).Value2 = 20.6042627779399
$ws.Cells.Item(11, 2).Value2 = 7.241016101368622
$ws.Cells.Item(11, 4).Value2 = 10.47482899352876
$ws.Cells.Item(11, 5).Value2 = 14.4411813493935
$ws.Cells.Item(11, 6).Value2 = 28.23594242244596
$ws.Cells.Item(11, 7).Value2 = 26.64400419352246
$ws.Cells.Item(11, 8).Value2 = 13.67937703808339
$ws.Cells.Item(11, 9).Value2 = 18.69353264686819
$ws.Cells.Item(11, 10).Value2 = 9.756647552340826
$ws.Cells.Item(11, 11).Value2 = 13.72878066605303
$ws.Cells.Item(11, 15).Value2 = 20.55383845358569
$ws.Cells.Item(12, 2).Value2 = 7.288422979949696
$ws.Cells.Item(12, 4).Value2 = 10.49895277905228
$ws.Cells.Item(12, 5).Value2 = 14.46828077785819
$ws.Cells.Item(12, 6).Value2 = 28.23181526328618
$ws.Cells.Item(12, 7).Value2 = 26.63363385725048
$ws.Cells.Item(12, 8).Value2 = 13.66714717763693
$ws.Cells.Item(12, 9).Value2 = 18.66611846449621
$ws.Cells.Item(12, 10).Value2 = 9.758027777913506
$ws.Cells.Item(12, 11).Value2 = 13.85631984419778
$ws.Cells.Item(12, 15).Value2 = 20.53552570394248
$ws.Cells.Item(13, 2).Value2 = 7.278237085304684
$ws.Cells.Item(13, 4).Value2 = 10.49375188854301
$ws.Cells.Item(13, 5).Value2 = 14.46243054396062
$ws.Cells.Item(13, 6).Value2 = 28.23265096962241
$ws.Cells.Item(13, 7).Value2 = 26.63579538590689
$ws.Cells.Item(13, 8).Value2 = 13.6697646116756
$ws.Cells.Item(13, 9).Value2 = 18.67199645655444
$ws.Cells.Item(13, 10).Value2 = 9.757720228156572
$ws.Cells.Item(13, 11).Value2 = 13.82895616478065
$ws.Cells.Item(13, 15).Value2 = 20.53943485562907
$ws.Cells.Item(14, 2).Value2 = 7.244927064041653
$ws.Cells.Item(14, 4).Value2 = 10.47681121278396
$ws.Cells.Item(14, 5).Value2 = 14.4434045501292
$ws.Cells.Item(14, 6).Value2 = 28.23557892958244
$ws.Cells.Item(14, 7).Value2 = 26.64311862970784
$ws.Cells.Item(14, 8).Value2 = 13.6783634520823
$ws.Cells.Item(14, 9).Value2 = 18.69126548347805
$ws.Cells.Item(14, 10).Value2 = 9.756756464176917
$ws.Cells.Item(14, 11).Value2 = 13.73932024078648
$ws.Cells.Item(14, 15).Value2 = 20.55231616391879
$ws.Cells.Item(15, 2).Value2 = 7.224454021917872
$ws.Cells.Item(15, 4).Value2 = 10.46645066911903
$ws.Cells.Item(15, 5).Value2 = 14.43179156731428
$ws.Cells.Item(15, 6).Value2 = 28.23752799895818
$ws.Cells.Item(15, 7).Value2 = 26.64781472664401
$ws.Cells.Item(15, 8).Value2 = 13.68367875442307
$ws.Cells.Item(15, 9).Value2 = 18.70314488777699
$ws.Cells.Item(15, 10).Value2 = 9.756196291358501
$ws.Cells.Item(15, 11).Value2 = 13.68411168743071
$ws.Cells.Item(15, 15).Value2 = 20.5603082571262
$ws.Cells.Item(16, 2).Value2 = 7.106092659483133
$ws.Cells.Item(16, 4).Value2 = 10.4073383452231
$ws.Cells.Item(16, 5).Value2 = 14.36588466617122
$ws.Cells.Item(16, 6).Value2 = 28.25108115755381
$ws.Cells.Item(16, 7).Value2 = 26.6779437065312
$ws.Cells.Item(16, 8).Value2 = 13.71487807499031
$ws.Cells.Item(16, 9).Value2 = 18.77239542978637
$ws.Cells.Item(16, 10).Value2 = 9.753449289017828
$ws.Cells.Item(16, 11).Value2 = 13.36312294637735
$ws.Cells.Item(16, 15).Value2 = 20.60766734984696
$ws.Cells.Item(17, 2).Value2 = 7.032588219330252
$ws.Cells.Item(17, 4).Value2 = 10.371320574359
$ws.Cells.Item(17, 5).Value2 = 14.32604141292919
$ws.Cells.Item(17, 6).Value2 = 28.26152113150635
$ws.Cells.Item(17, 7).Value2 = 26.69928962390651
$ws.Cells.Item(17, 8).Value2 = 13.73467738151262
$ws.Cells.Item(17, 9).Value2 = 18.81592738204233
$ws.Cells.Item(17, 10).Value2 = 9.752178224046707
$ws.Cells.Item(17, 11).Value2 = 13.16214785862148
$ws.Cells.Item(17, 15).Value2 = 20.63810955746945
$ws.Cells.Item(18, 2).Value2 = 6.989995645216101
$ws.Cells.Item(18, 4).Value2 = 10.35069991299736
$ws.Cells.Item(18, 5).Value2 = 14.30334688142353
$ws.Cells.Item(18, 6).Value2 = 28.26830633370574
$ws.Cells.Item(18, 7).Value2 = 26.7126158345276
$ws.Cells.Item(18, 8).Value2 = 13.74630767396495
$ws.Cells.Item(18, 9).Value2 = 18.8413517040985
$ws.Cells.Item(18, 10).Value2 = 9.751599666452599
$ws.Cells.Item(18, 11).Value2 = 13.04508404632249
$ws.Cells.Item(18, 15).Value2 = 20.6561284409367
$ws.Cells.Item(19, 2).Value2 = 6.975521805419145
$ws.Cells.Item(19, 4).Value2 = 10.34373508639037
$ws.Cells.Item(19, 5).Value2 = 14.29570163616437
$ws.Cells.Item(19, 6).Value2 = 28.27073764384649
$ws.Cells.Item(19, 7).Value2 = 26.71730761724463
$ws.Cells.Item(19, 8).Value2 = 13.75028708733022
$ws.Cells.Item(19, 9).Value2 = 18.85002625856114
$ws.Cells.Item(19, 10).Value2 = 9.751430003310697
$ws.Cells.Item(19, 11).Value2 = 13.00519842875789
$ws.Cells.Item(19, 15).Value2 = 20.66231669846677
$ws.Cells.Item(20, 2).Value2 = 7.040445784009894
$ws.Cells.Item(20, 4).Value2 = 10.37514493556426
$ws.Cells.Item(20, 5).Value2 = 14.33025992657783
$ws.Cells.Item(20, 6).Value2 = 28.26032899142697
$ws.Cells.Item(20, 7).Value2 = 26.6969087076334
$ws.Cells.Item(20, 8).Value2 = 13.73254463632784
$ws.Cells.Item(20, 9).Value2 = 18.81125340179688
$ws.Cells.Item(20, 10).Value2 = 9.752297752539908
$ws.Cells.Item(20, 11).Value2 = 13.18369450926079
$ws.Cells.Item(20, 15).Value2 = 20.63481619376189
$ws.Cells.Item(21, 2).Value2 = 7.254725617897736
$ws.Cells.Item(21, 4).Value2 = 10.48178377352417
$ws.Cells.Item(21, 5).Value2 = 14.4489844380296
$ws.Cells.Item(21, 6).Value2 = 28.23468648526688
$ws.Cells.Item(21, 7).Value2 = 26.640923751646
$ws.Cells.Item(21, 8).Value2 = 13.6758277052802
$ws.Cells.Item(21, 9).Value2 = 18.68558974838346
$ws.Cells.Item(21, 10).Value2 = 9.757033261679501
$ws.Cells.Item(21, 11).Value2 = 13.76571194081932
$ws.Cells.Item(21, 15).Value2 = 20.54851136539603
$ws.Cells.Item(22, 2).Value2 = 7.391686491217044
$ws.Cells.Item(22, 4).Value2 = 10.55221429322346
$ws.Cells.Item(22, 5).Value2 = 14.52842851204531
$ws.Cells.Item(22, 6).Value2 = 28.22489053463476
$ws.Cells.Item(22, 7).Value2 = 26.61374225184043
$ws.Cells.Item(22, 8).Value2 = 13.64091977858878
$ws.Cells.Item(22, 9).Value2 = 18.60688938150007
$ws.Cells.Item(22, 10).Value2 = 9.761479023249301
$ws.Cells.Item(22, 11).Value2 = 14.13254971969801
$ws.Cells.Item(22, 15).Value2 = 20.49666491868753
$ws.Cells.Item(23, 2).Value2 = 7.318882871678708
$ws.Cells.Item(23, 4).Value2 = 10.51456256788167
$ws.Cells.Item(23, 5).Value2 = 14.48586468759151
$ws.Cells.Item(23, 6).Value2 = 28.22948121961317
$ws.Cells.Item(23, 7).Value2 = 26.62738556734108
$ws.Cells.Item(23, 8).Value2 = 13.65935303885494
$ws.Cells.Item(23, 9).Value2 = 18.64857996521359
$ws.Cells.Item(23, 10).Value2 = 9.758983018255691
$ws.Cells.Item(23, 11).Value2 = 13.93802138126363
$ws.Cells.Item(23, 15).Value2 = 20.52391810508396
$ws.Cells.Item(24, 2).Value2 = 7.036894415814234
$ws.Cells.Item(24, 4).Value2 = 10.37341567133357
$ws.Cells.Item(24, 5).Value2 = 14.32835207515437
$ws.Cells.Item(24, 6).Value2 = 28.26086551875326
$ws.Cells.Item(24, 7).Value2 = 26.69798183705822
$ws.Cells.Item(24, 8).Value2 = 13.73350807978601
$ws.Cells.Item(24, 9).Value2 = 18.81336527085983
$ws.Cells.Item(24, 10).Value2 = 9.75224323956081
$ws.Cells.Item(24, 11).Value2 = 13.17395799920415
$ws.Cells.Item(24, 15).Value2 = 20.63630351253056
$ws.Cells.Item(25, 2).Value2 = 6.720892700492636
$ws.Cells.Item(25, 4).Value2 = 10.22489583345978
$ws.Cells.Item(25, 5).Value2 = 14.16701599230873
$ws.Cells.Item(25, 6).Value2 = 28.32453876374265
$ws.Cells.Item(25, 7).Value2 = 26.81416710093865
$ws.Cells.Item(25, 8).Value2 = 13.82275689951304
$ws.Cells.Item(25, 9).Value2 = 19.00587226443685
$ws.Cells.Item(25, 10).Value2 = 9.750807639949146
$ws.Cells.Item(25, 11).Value2 = 12.2940530901943
$ws.Cells.Item(25, 15).Value2 = 20.77698789190762
